$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Sheet1" to "end_point"
$ws.Name = "end_point"

# Update cell C13: "Get - provide all server (conditional)" -> new text
$ws.Range("C13").Value = "Get - check (conditional) server's ip exists or not"
